$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out the old data (rows 3 and 4 are no longer used)
$ws.Range("A1:D4").ClearContents()

# Header row
$ws.Range("A1").Value = "y"
$ws.Range("B1").Value = "m"
$ws.Range("C1").Value = "d"
$ws.Range("D1").Value = "yToD"
$ws.Range("E1").Value = "yToM"
$ws.Range("F1").Value = "days"

# Data row
$ws.Range("A2").Value = 25
$ws.Range("B2").Value = 6
$ws.Range("C2").Value = 15
$ws.Range("D2").Value = 9125
$ws.Range("E2").Value = 180
$ws.Range("F2").Value = 9320

$ws.Range("G5").Select()
